$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B19").Value = "발표 녹화, 영상촬영"
$ws.Range("D19").Value = "ppt제작"
$ws.Range("F19").Value = "영상편집"
$ws.Range("H19").Value = "영상촬영"
$ws.Range("J19").Value = "자연어 정규표현식 공부"
